$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-9 (columns A-T), reflecting the
# expanded NATMI Il12b -> Il12rb1 signalling table across the FAPs/M2
# sending clusters and ECs/FAPs/M2/sCs target clusters.
$rows = @(
    @("FAPs","Il12b","Il12rb1","ECs",1,0.3333333333333333,0.091959,0.275877,0.06282071922171414,0.06282071922171413,2,0.6666666666666666,0.162136,0.486408,0.0598694021675715,0.0598694021675715,0.014909864424,0.134188779816,0.003761038903540894,0.003761038903540893),
    @("FAPs","Il12b","Il12rb1","FAPs",1,0.3333333333333333,0.091959,0.275877,0.06282071922171414,0.06282071922171413,3,1,1.694965,5.084895,0.6258729785177741,0.6258729785177741,0.155867286435,1.402805577915,0.03931779065192301,0.03931779065192301),
    @("FAPs","Il12b","Il12rb1","M2",1,0.3333333333333333,0.091959,0.275877,0.06282071922171414,0.06282071922171413,3,1,0.6978413333333333,2.093524,0.2576808570242738,0.2576808570242737,0.064172791172,0.577555120548,0.01618769676793257,0.01618769676793256),
    @("FAPs","Il12b","Il12rb1","sCs",1,0.3333333333333333,0.091959,0.275877,0.06282071922171414,0.06282071922171413,2,0.6666666666666666,0.153219,0.459657,0.05657676229038053,0.05657676229038053,0.014089866021,0.126808794189,0.00355419289831766,0.003554192898317659),
    @("M2","Il12b","Il12rb1","ECs",3,1,1.371873333333334,4.115620000000001,0.9371792807782859,0.9371792807782858,2,0.6666666666666666,0.162136,0.486408,0.0598694021675715,0.0598694021675715,0.2224300547733334,2.00187049296,0.05610836326403061,0.0561083632640306),
    @("M2","Il12b","Il12rb1","FAPs",3,1,1.371873333333334,4.115620000000001,0.9371792807782859,0.9371792807782858,3,1,1.694965,5.084895,0.6258729785177741,0.6258729785177741,2.325277284433334,20.92749555990001,0.5865551878658511,0.5865551878658511),
    @("M2","Il12b","Il12rb1","M2",3,1,1.371873333333334,4.115620000000001,0.9371792807782859,0.9371792807782858,3,1,0.6978413333333333,2.093524,0.2576808570242738,0.2576808570242737,0.9573499160977779,8.61614924488,0.2414931602563412,0.2414931602563412),
    @("M2","Il12b","Il12rb1","sCs",3,1,1.371873333333334,4.115620000000001,0.9371792807782859,0.9371792807782858,2,0.6666666666666666,0.153219,0.459657,0.05657676229038053,0.05657676229038053,0.21019706026,1.89177354234,0.05302256939206287,0.05302256939206287)
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowData = $rows[$r]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
    }
}

Write-Output "Populated rows 2-9 with updated Il12b/Il12rb1 NATMI results"
